$wb = $excel.ActiveWorkbook

# ARM sheet - row 122 updates
$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Cells.Item(122, 8).Value = 6685.0557
$wsARM.Cells.Item(122, 9).Value = 7987.7856
$wsARM.Cells.Item(122, 10).Value = 2125.5
$wsARM.Cells.Item(122, 11).Value = 23963.3568
$wsARM.Cells.Item(122, 12).Value = 6376.5
$wsARM.Cells.Item(122, 13).Value = -21513.3568
$wsARM.Cells.Item(122, 14).Value = -11276.5

# CUL sheet - new H:N values for rows 120-141 (row 135 already populated)
$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Cells.Item(120, 8).Value = 48006
$wsCUL.Cells.Item(120, 9).Value = 100030
$wsCUL.Cells.Item(120, 10).Value = 35000
$wsCUL.Cells.Item(120, 11).Value = 300090
$wsCUL.Cells.Item(120, 12).Value = 105000
$wsCUL.Cells.Item(120, 13).Value = -295252
$wsCUL.Cells.Item(120, 14).Value = -114676
$wsCUL.Cells.Item(121, 8).Value = 195
$wsCUL.Cells.Item(121, 9).Value = 195
$wsCUL.Cells.Item(121, 10).Value = 0
$wsCUL.Cells.Item(121, 11).Value = 585
$wsCUL.Cells.Item(121, 12).Value = 0
$wsCUL.Cells.Item(121, 13).Value = 725
$wsCUL.Cells.Item(122, 8).Value = 790.8
$wsCUL.Cells.Item(122, 9).Value = 472.85715
$wsCUL.Cells.Item(122, 10).Value = 1532.6666
$wsCUL.Cells.Item(122, 11).Value = 4255.71435
$wsCUL.Cells.Item(122, 12).Value = 13793.9994
$wsCUL.Cells.Item(122, 13).Value = -1805.71435
$wsCUL.Cells.Item(122, 14).Value = -18693.9994
$wsCUL.Cells.Item(123, 8).Value = 3000
$wsCUL.Cells.Item(123, 9).Value = 2000
$wsCUL.Cells.Item(123, 10).Value = 4000
$wsCUL.Cells.Item(123, 11).Value = 6000
$wsCUL.Cells.Item(123, 12).Value = 12000
$wsCUL.Cells.Item(123, 13).Value = -3550
$wsCUL.Cells.Item(123, 14).Value = -16900
$wsCUL.Cells.Item(124, 8).Value = 1033.3334
$wsCUL.Cells.Item(124, 9).Value = 750
$wsCUL.Cells.Item(124, 10).Value = 1600
$wsCUL.Cells.Item(124, 11).Value = 2250
$wsCUL.Cells.Item(124, 12).Value = 4800
$wsCUL.Cells.Item(124, 13).Value = 2660
$wsCUL.Cells.Item(124, 14).Value = -14620
$wsCUL.Cells.Item(125, 8).Value = 2469.7742
$wsCUL.Cells.Item(125, 9).Value = 1283.3334
$wsCUL.Cells.Item(125, 10).Value = 2596.8928
$wsCUL.Cells.Item(125, 11).Value = 3850.0002
$wsCUL.Cells.Item(125, 12).Value = 7790.678400000001
$wsCUL.Cells.Item(125, 13).Value = 1069.9998
$wsCUL.Cells.Item(125, 14).Value = -17630.6784
$wsCUL.Cells.Item(126, 8).Value = 250002510
$wsCUL.Cells.Item(126, 9).Value = 1030
$wsCUL.Cells.Item(126, 10).Value = 333336320
$wsCUL.Cells.Item(126, 11).Value = 3090
$wsCUL.Cells.Item(126, 12).Value = 1000008960
$wsCUL.Cells.Item(126, 13).Value = 1850
$wsCUL.Cells.Item(126, 14).Value = -1000018840
$wsCUL.Cells.Item(127, 8).Value = 2750
$wsCUL.Cells.Item(127, 9).Value = 0
$wsCUL.Cells.Item(127, 10).Value = 2750
$wsCUL.Cells.Item(127, 11).Value = 0
$wsCUL.Cells.Item(127, 12).Value = 8250
$wsCUL.Cells.Item(127, 14).Value = -18170
$wsCUL.Cells.Item(128, 8).Value = 95666.336
$wsCUL.Cells.Item(128, 9).Value = 95666.336
$wsCUL.Cells.Item(128, 10).Value = 0
$wsCUL.Cells.Item(128, 11).Value = 286999.008
$wsCUL.Cells.Item(128, 12).Value = 0
$wsCUL.Cells.Item(128, 13).Value = -282019.008
$wsCUL.Cells.Item(129, 8).Value = 1106.2307
$wsCUL.Cells.Item(129, 9).Value = 562.5
$wsCUL.Cells.Item(129, 10).Value = 1347.8889
$wsCUL.Cells.Item(129, 11).Value = 1687.5
$wsCUL.Cells.Item(129, 12).Value = 4043.6667
$wsCUL.Cells.Item(129, 13).Value = 3312.5
$wsCUL.Cells.Item(129, 14).Value = -14043.6667
$wsCUL.Cells.Item(130, 8).Value = 998
$wsCUL.Cells.Item(130, 9).Value = 0
$wsCUL.Cells.Item(130, 10).Value = 998
$wsCUL.Cells.Item(130, 11).Value = 0
$wsCUL.Cells.Item(130, 12).Value = 2994
$wsCUL.Cells.Item(130, 14).Value = -13034
$wsCUL.Cells.Item(131, 8).Value = 2912.0146
$wsCUL.Cells.Item(131, 9).Value = 315
$wsCUL.Cells.Item(131, 10).Value = 2990.7122
$wsCUL.Cells.Item(131, 11).Value = 945
$wsCUL.Cells.Item(131, 12).Value = 8972.1366
$wsCUL.Cells.Item(131, 13).Value = 4095
$wsCUL.Cells.Item(131, 14).Value = -19052.1366
$wsCUL.Cells.Item(132, 8).Value = 961.7222
$wsCUL.Cells.Item(132, 9).Value = 931.1
$wsCUL.Cells.Item(132, 10).Value = 1000
$wsCUL.Cells.Item(132, 11).Value = 8379.9
$wsCUL.Cells.Item(132, 12).Value = 9000
$wsCUL.Cells.Item(132, 13).Value = -5849.9
$wsCUL.Cells.Item(132, 14).Value = -14060
$wsCUL.Cells.Item(133, 8).Value = 9666.666999999999
$wsCUL.Cells.Item(133, 9).Value = 8500
$wsCUL.Cells.Item(133, 10).Value = 12000
$wsCUL.Cells.Item(133, 11).Value = 25500
$wsCUL.Cells.Item(133, 12).Value = 36000
$wsCUL.Cells.Item(133, 13).Value = -20440
$wsCUL.Cells.Item(133, 14).Value = -46120
$wsCUL.Cells.Item(134, 8).Value = 2759.8
$wsCUL.Cells.Item(134, 9).Value = 2949.75
$wsCUL.Cells.Item(134, 10).Value = 2000
$wsCUL.Cells.Item(134, 11).Value = 8849.25
$wsCUL.Cells.Item(134, 12).Value = 6000
$wsCUL.Cells.Item(134, 13).Value = -3779.25
$wsCUL.Cells.Item(134, 14).Value = -16140
$wsCUL.Cells.Item(136, 8).Value = 5186.647
$wsCUL.Cells.Item(136, 9).Value = 1140
$wsCUL.Cells.Item(136, 10).Value = 6431.769
$wsCUL.Cells.Item(136, 11).Value = 3420
$wsCUL.Cells.Item(136, 12).Value = 19295.307
$wsCUL.Cells.Item(136, 13).Value = 1680
$wsCUL.Cells.Item(136, 14).Value = -29495.307
$wsCUL.Cells.Item(137, 8).Value = 5614495
$wsCUL.Cells.Item(137, 9).Value = 9092395
$wsCUL.Cells.Item(137, 10).Value = 149223.72
$wsCUL.Cells.Item(137, 11).Value = 27277185
$wsCUL.Cells.Item(137, 12).Value = 447671.16
$wsCUL.Cells.Item(137, 13).Value = -27272085
$wsCUL.Cells.Item(137, 14).Value = -457871.16
$wsCUL.Cells.Item(138, 8).Value = 874.9286
$wsCUL.Cells.Item(138, 9).Value = 849.9231
$wsCUL.Cells.Item(138, 10).Value = 1200
$wsCUL.Cells.Item(138, 11).Value = 2549.7693
$wsCUL.Cells.Item(138, 12).Value = 3600
$wsCUL.Cells.Item(138, 13).Value = 2590.2307
$wsCUL.Cells.Item(138, 14).Value = -13880
$wsCUL.Cells.Item(139, 8).Value = 38463228
$wsCUL.Cells.Item(139, 9).Value = 41668080
$wsCUL.Cells.Item(139, 10).Value = 4977
$wsCUL.Cells.Item(139, 11).Value = 125004240
$wsCUL.Cells.Item(139, 12).Value = 14931
$wsCUL.Cells.Item(139, 13).Value = -124999100
$wsCUL.Cells.Item(139, 14).Value = -25211
$wsCUL.Cells.Item(140, 8).Value = 9507.916999999999
$wsCUL.Cells.Item(140, 9).Value = 11525.556
$wsCUL.Cells.Item(140, 10).Value = 3455
$wsCUL.Cells.Item(140, 11).Value = 34576.66800000001
$wsCUL.Cells.Item(140, 12).Value = 10365
$wsCUL.Cells.Item(140, 13).Value = -29396.66800000001
$wsCUL.Cells.Item(140, 14).Value = -20725
$wsCUL.Cells.Item(141, 8).Value = 3879
$wsCUL.Cells.Item(141, 9).Value = 4016.9
$wsCUL.Cells.Item(141, 10).Value = 2500
$wsCUL.Cells.Item(141, 11).Value = 12050.7
$wsCUL.Cells.Item(141, 12).Value = 7500
$wsCUL.Cells.Item(141, 13).Value = -6870.700000000001
$wsCUL.Cells.Item(141, 14).Value = -17860

# GSM sheet - row 122 updates, and new H:N values for rows 125, 127-141
$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Cells.Item(122, 8).Value = 1012419.56
$wsGSM.Cells.Item(122, 9).Value = 1012419.56
$wsGSM.Cells.Item(122, 11).Value = 3037258.68
$wsGSM.Cells.Item(122, 13).Value = -3034808.68
$wsGSM.Cells.Item(125, 8).Value = 0
$wsGSM.Cells.Item(125, 9).Value = 0
$wsGSM.Cells.Item(125, 10).Value = 0
$wsGSM.Cells.Item(125, 11).Value = 0
$wsGSM.Cells.Item(125, 12).Value = 0
$wsGSM.Cells.Item(126, 8).Value = 2143.4
$wsGSM.Cells.Item(126, 9).Value = 1813.9412
$wsGSM.Cells.Item(126, 10).Value = 2386.913
$wsGSM.Cells.Item(126, 11).Value = 5441.8236
$wsGSM.Cells.Item(126, 12).Value = 7160.739
$wsGSM.Cells.Item(126, 13).Value = -2971.8236
$wsGSM.Cells.Item(126, 14).Value = -12100.739
$wsGSM.Cells.Item(127, 8).Value = 0
$wsGSM.Cells.Item(127, 9).Value = 0
$wsGSM.Cells.Item(127, 10).Value = 0
$wsGSM.Cells.Item(127, 11).Value = 0
$wsGSM.Cells.Item(127, 12).Value = 0
$wsGSM.Cells.Item(128, 8).Value = 54500
$wsGSM.Cells.Item(128, 9).Value = 0
$wsGSM.Cells.Item(128, 10).Value = 54500
$wsGSM.Cells.Item(128, 11).Value = 0
$wsGSM.Cells.Item(128, 12).Value = 54500
$wsGSM.Cells.Item(128, 14).Value = -64460
$wsGSM.Cells.Item(129, 8).Value = 49999.8
$wsGSM.Cells.Item(129, 9).Value = 0
$wsGSM.Cells.Item(129, 10).Value = 49999.8
$wsGSM.Cells.Item(129, 11).Value = 0
$wsGSM.Cells.Item(129, 12).Value = 49999.8
$wsGSM.Cells.Item(129, 14).Value = -59999.8
$wsGSM.Cells.Item(130, 8).Value = 49800
$wsGSM.Cells.Item(130, 9).Value = 0
$wsGSM.Cells.Item(130, 10).Value = 49800
$wsGSM.Cells.Item(130, 11).Value = 0
$wsGSM.Cells.Item(130, 12).Value = 49800
$wsGSM.Cells.Item(130, 14).Value = -59840
$wsGSM.Cells.Item(131, 8).Value = 35000
$wsGSM.Cells.Item(131, 9).Value = 0
$wsGSM.Cells.Item(131, 10).Value = 35000
$wsGSM.Cells.Item(131, 11).Value = 0
$wsGSM.Cells.Item(131, 12).Value = 35000
$wsGSM.Cells.Item(131, 14).Value = -45080
$wsGSM.Cells.Item(132, 8).Value = 3309.5715
$wsGSM.Cells.Item(132, 9).Value = 3030.5386
$wsGSM.Cells.Item(132, 10).Value = 3763
$wsGSM.Cells.Item(132, 11).Value = 9091.6158
$wsGSM.Cells.Item(132, 12).Value = 11289
$wsGSM.Cells.Item(132, 13).Value = -6561.6158
$wsGSM.Cells.Item(132, 14).Value = -16349
$wsGSM.Cells.Item(133, 8).Value = 0
$wsGSM.Cells.Item(133, 9).Value = 0
$wsGSM.Cells.Item(133, 10).Value = 0
$wsGSM.Cells.Item(133, 11).Value = 0
$wsGSM.Cells.Item(133, 12).Value = 0
$wsGSM.Cells.Item(134, 8).Value = 0
$wsGSM.Cells.Item(134, 9).Value = 0
$wsGSM.Cells.Item(134, 10).Value = 0
$wsGSM.Cells.Item(134, 11).Value = 0
$wsGSM.Cells.Item(134, 12).Value = 0
$wsGSM.Cells.Item(135, 8).Value = 1000000000
$wsGSM.Cells.Item(135, 9).Value = 0
$wsGSM.Cells.Item(135, 10).Value = 1000000000
$wsGSM.Cells.Item(135, 11).Value = 0
$wsGSM.Cells.Item(135, 12).Value = 1000000000
$wsGSM.Cells.Item(135, 14).Value = -1000010140
$wsGSM.Cells.Item(136, 8).Value = 40326
$wsGSM.Cells.Item(136, 9).Value = 0
$wsGSM.Cells.Item(136, 10).Value = 40326
$wsGSM.Cells.Item(136, 11).Value = 0
$wsGSM.Cells.Item(136, 12).Value = 120978
$wsGSM.Cells.Item(136, 14).Value = -126078
$wsGSM.Cells.Item(137, 8).Value = 0
$wsGSM.Cells.Item(137, 9).Value = 0
$wsGSM.Cells.Item(137, 10).Value = 0
$wsGSM.Cells.Item(137, 11).Value = 0
$wsGSM.Cells.Item(137, 12).Value = 0
$wsGSM.Cells.Item(138, 8).Value = 0
$wsGSM.Cells.Item(138, 9).Value = 0
$wsGSM.Cells.Item(138, 10).Value = 0
$wsGSM.Cells.Item(138, 11).Value = 0
$wsGSM.Cells.Item(138, 12).Value = 0
$wsGSM.Cells.Item(139, 8).Value = 0
$wsGSM.Cells.Item(139, 9).Value = 0
$wsGSM.Cells.Item(139, 10).Value = 0
$wsGSM.Cells.Item(139, 11).Value = 0
$wsGSM.Cells.Item(139, 12).Value = 0
$wsGSM.Cells.Item(140, 8).Value = 51832.5
$wsGSM.Cells.Item(140, 9).Value = 0
$wsGSM.Cells.Item(140, 10).Value = 51832.5
$wsGSM.Cells.Item(140, 11).Value = 0
$wsGSM.Cells.Item(140, 12).Value = 51832.5
$wsGSM.Cells.Item(140, 14).Value = -62192.5
$wsGSM.Cells.Item(141, 8).Value = 68400
$wsGSM.Cells.Item(141, 9).Value = 0
$wsGSM.Cells.Item(141, 10).Value = 68400
$wsGSM.Cells.Item(141, 11).Value = 0
$wsGSM.Cells.Item(141, 12).Value = 68400
$wsGSM.Cells.Item(141, 14).Value = -78760
